$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed by Excel
# (e.g. "91%") must be pre-formatted as Text so the literal string is kept.

$ws.Range("E2").Value = "2026-02-20 07:07:57"
$ws.Range("E3").Value = "2026-02-20 07:08:00"
$ws.Range("E4").Value = "2026-02-20 07:08:02"
$ws.Range("E5").Value = "2026-02-20 07:08:05"
$ws.Range("E6").Value = "2026-02-20 07:08:08"
$ws.Range("E7").Value = "2026-02-20 07:08:10"
$ws.Range("E8").Value = "2026-02-20 07:08:13"
$ws.Range("E9").Value = "2026-02-20 07:08:15"
$ws.Range("E10").Value = "2026-02-20 07:08:18"
$ws.Range("E11").Value = "2026-02-20 07:08:20"
$ws.Range("E12").Value = "2026-02-20 07:08:23"
$ws.Range("E13").Value = "2026-02-20 07:08:25"
$ws.Range("E14").Value = "2026-02-20 07:08:28"
$ws.Range("E15").Value = "2026-02-20 07:08:31"
$ws.Range("E16").Value = "2026-02-20 07:08:33"
$ws.Range("E17").Value = "2026-02-20 07:08:36"
$ws.Range("E18").Value = "2026-02-20 07:08:38"
$ws.Range("E19").Value = "2026-02-20 07:08:41"
$ws.Range("E20").Value = "2026-02-20 07:08:44"
$ws.Range("E21").Value = "2026-02-20 07:08:47"
$ws.Range("E22").Value = "2026-02-20 07:08:49"
$ws.Range("E23").Value = "2026-02-20 07:08:52"
$ws.Range("E24").Value = "2026-02-20 07:08:55"
$ws.Range("E25").Value = "2026-02-20 07:08:58"
$ws.Range("E26").Value = "2026-02-20 07:09:00"
$ws.Range("E27").Value = "2026-02-20 07:09:03"
$ws.Range("E28").Value = "2026-02-20 07:09:06"
$ws.Range("J28").Value = "1020.2 hPa"
$ws.Range("N28").Value = "0.0 °C 6:33 TU"
$ws.Range("O28").Value = "2.3 °C"
$ws.Range("E29").Value = "2026-02-20 07:09:09"
$ws.Range("H29").NumberFormat = "@"
$ws.Range("H29").Value = "91%"
$ws.Range("I29").Value = "0.1 mm"
$ws.Range("O29").Value = "3.6 °C"
$ws.Range("E30").Value = "2026-02-20 07:09:11"
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = "65%"
$ws.Range("J30").Value = "1018.7 hPa"
$ws.Range("K30").Value = "0.0 MJ/m2"
$ws.Range("O30").Value = "8.8 °C"
$ws.Range("E31").Value = "2026-02-20 07:09:14"
$ws.Range("H31").NumberFormat = "@"
$ws.Range("H31").Value = "57%"
$ws.Range("J31").Value = "1017.1 hPa"
$ws.Range("K31").Value = "0.0 MJ/m2"
$ws.Range("E32").Value = "2026-02-20 07:09:17"
$ws.Range("H32").NumberFormat = "@"
$ws.Range("H32").Value = "91%"
$ws.Range("N32").Value = "1.1 °C 6:40 TU"
$ws.Range("E33").Value = "2026-02-20 07:09:19"
$ws.Range("H33").NumberFormat = "@"
$ws.Range("H33").Value = "47%"
$ws.Range("J33").Value = "1019.7 hPa"
$ws.Range("O33").Value = "4.3 °C"
$ws.Range("E34").Value = "2026-02-20 07:09:22"
$ws.Range("H34").NumberFormat = "@"
$ws.Range("H34").Value = "60%"
$ws.Range("E35").Value = "2026-02-20 07:09:25"
$ws.Range("J35").Value = "1024.4 hPa"
$ws.Range("E36").Value = "2026-02-20 07:09:27"
$ws.Range("J36").Value = "1019.0 hPa"
$ws.Range("N36").Value = "12.2 °C 6:58 TU"
$ws.Range("O36").Value = "13.4 °C"
$ws.Range("E37").Value = "2026-02-20 07:09:30"
$ws.Range("H37").NumberFormat = "@"
$ws.Range("H37").Value = "69%"
$ws.Range("J37").Value = "1021.9 hPa"
$ws.Range("N37").Value = "-1.1 °C 6:57 TU"
$ws.Range("O37").Value = "1.5 °C"
$ws.Range("E38").Value = "2026-02-20 07:09:33"
$ws.Range("K38").Value = "0.0 MJ/m2"
$ws.Range("L38").Value = "14.0 km/h - 300º 6:46 TU"
$ws.Range("O38").Value = "4.0 °C"
$ws.Range("E39").Value = "2026-02-20 07:09:35"
$ws.Range("H39").NumberFormat = "@"
$ws.Range("H39").Value = "70%"
$ws.Range("M39").Value = "-3.8 °C 6:40 TU"
$ws.Range("O39").Value = "-5.7 °C"
$ws.Range("E40").Value = "2026-02-20 07:09:37"
$ws.Range("H40").NumberFormat = "@"
$ws.Range("H40").Value = "46%"
$ws.Range("J40").Value = "1021.2 hPa"
$ws.Range("M40").Value = "9.5 °C 6:45 TU"
$ws.Range("O40").Value = "7.6 °C"
$ws.Range("E41").Value = "2026-02-20 07:09:40"
$ws.Range("H41").NumberFormat = "@"
$ws.Range("H41").Value = "52%"
$ws.Range("J41").Value = "1020.6 hPa"
$ws.Range("N41").Value = "9.9 °C 6:58 TU"
$ws.Range("O41").Value = "10.7 °C"
$ws.Range("E42").Value = "2026-02-20 07:09:43"
$ws.Range("H42").NumberFormat = "@"
$ws.Range("H42").Value = "92%"
$ws.Range("E43").Value = "2026-02-20 07:09:45"
$ws.Range("H43").NumberFormat = "@"
$ws.Range("H43").Value = "86%"
$ws.Range("K43").Value = "0.0 MJ/m2"
$ws.Range("N43").Value = "-0.3 °C 6:45 TU"
$ws.Range("O43").Value = "1.1 °C"
$ws.Range("E44").Value = "2026-02-20 07:09:48"
$ws.Range("I44").Value = "4.8 mm"
$ws.Range("E45").Value = "2026-02-20 07:09:51"
$ws.Range("J45").Value = "1027.9 hPa"
$ws.Range("N45").Value = "1.6 °C 6:51 TU"
$ws.Range("E46").Value = "2026-02-20 07:09:53"
$ws.Range("J46").Value = "1024.2 hPa"
$ws.Range("N46").Value = "8.6 °C 6:30 TU"
$ws.Range("O46").Value = "9.5 °C"
